$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.WindowState = -4137
$win.Width = 15345
$win.Height = 4455
$win.UsableWidth = 15345
$win.UsableHeight = 4455
